# Adds the person-details row submitted from the order page (row 3):
# Customer ID, Name, Phone, Address -> text values, Phone/Address left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new row's cells as Text first so the numeric-looking
# customer id is stored verbatim (not coerced into a Number), and so the
# blank Phone/Address cells are still written out instead of being
# dropped as "no value".
$ws.Range("A3:D3").NumberFormat = "@"

$ws.Range("A3").Value = "20250308122656"
$ws.Range("B3").Value = "a"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
